# Weekly driver report update for 2025-04-21
# Updates Critical Minutes / Good Roaming Calculation figures on the
# "Bad Drivers" table, and the Total Samples figure for one of the
# "Good Drivers" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers table (rows 3-5), Critical Minutes column (C)
$ws.Range("C3").Value = 10982
$ws.Range("C4").Value = 3846
$ws.Range("C5").Value = 200

# Good Roaming Calculation (%) for row 5 updated alongside the Critical
# Minutes change above
$ws.Range("D5").Value = 98.7

# Totals row recalculated Critical Minutes total
$ws.Range("C6").Value = 15028

# Good Drivers table: Total Samples for the
# "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8" row
$ws.Range("B16").Value = 331283
